# issue #5: stock data output to json file
#
# Adds a "property_category" column to the 股票 (stock) sheet, populated
# with the constant value "stock" for every data row. The new column is
# inserted right after the "total" column and before the "date" column,
# so the existing date / legislator_name / legislator_id columns shift
# one column to the right (H->I, I->J, J->K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column at H; this shifts the old H:J (date,
# legislator_name, legislator_id) columns to I:K and carries their
# values/styles along with them.
$ws.Columns("H:H").Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "property_category"

# Populate the new column for every existing data row with "stock".
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
$ws.Range("H5").Value = "stock"
